$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- C31 becomes a literal text "0" (percent calc undefined), matching D14's existing text-"0" style ---
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Weekly CompStat crime-count table updates (rows 14-31, 33; columns C..N) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 17
$ws.Range("K14").Value = -32
$ws.Range("L14").Value = -52.777777777777
$ws.Range("M14").Value = -55.263157894736
$ws.Range("N14").Value = -85.217391304347
# Row 15
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -28.571428571428
$ws.Range("F15").Value = 21
$ws.Range("H15").Value = 16.666666666666
$ws.Range("I15").Value = 125
$ws.Range("J15").Value = 104
$ws.Range("K15").Value = 20.192307692307
$ws.Range("L15").Value = 19.047619047619
$ws.Range("M15").Value = 56.25
$ws.Range("N15").Value = -54.545454545454
# Row 16
$ws.Range("C16").Value = 36
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 142
$ws.Range("G16").Value = 164
$ws.Range("H16").Value = -13.414634146341
$ws.Range("I16").Value = 763
$ws.Range("J16").Value = 897
$ws.Range("K16").Value = -14.938684503901
$ws.Range("L16").Value = -11.484918793503
$ws.Range("M16").Value = -42.501883948756
$ws.Range("N16").Value = -88.827061063113
# Row 17
$ws.Range("C17").Value = 90
$ws.Range("D17").Value = 76
$ws.Range("E17").Value = 18.421052631578
$ws.Range("F17").Value = 319
$ws.Range("G17").Value = 343
$ws.Range("H17").Value = -6.997084548104
$ws.Range("I17").Value = 1856
$ws.Range("J17").Value = 1786
$ws.Range("K17").Value = 3.919372900335
$ws.Range("L17").Value = 8.411214953271
$ws.Range("M17").Value = 61.531766753698
$ws.Range("N17").Value = -42.627511591962
# Row 18
$ws.Range("C18").Value = 27
$ws.Range("D18").Value = 28
$ws.Range("E18").Value = -3.571428571428
$ws.Range("F18").Value = 90
$ws.Range("G18").Value = 110
$ws.Range("H18").Value = -18.181818181818
$ws.Range("I18").Value = 704
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = -6.133333333333
$ws.Range("L18").Value = -19.542857142857
$ws.Range("M18").Value = -56.13707165109
$ws.Range("N18").Value = -92.195987141115
# Row 19
$ws.Range("C19").Value = 97
$ws.Range("D19").Value = 106
$ws.Range("E19").Value = -8.490566037735
$ws.Range("F19").Value = 417
$ws.Range("G19").Value = 481
$ws.Range("H19").Value = -13.305613305613
$ws.Range("I19").Value = 2590
$ws.Range("J19").Value = 3010
$ws.Range("K19").Value = -13.953488372093
$ws.Range("L19").Value = -20.770877944325
$ws.Range("M19").Value = 1.171875
$ws.Range("N19").Value = -38.082715754243
# Row 20
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 50
$ws.Range("E20").Value = -26
$ws.Range("F20").Value = 168
$ws.Range("G20").Value = 176
$ws.Range("H20").Value = -4.545454545454
$ws.Range("I20").Value = 833
$ws.Range("J20").Value = 973
$ws.Range("K20").Value = -14.388489208633
$ws.Range("L20").Value = -4.362801377726
$ws.Range("M20").Value = -13.946280991735
$ws.Range("N20").Value = -92.929292929292
# Row 21
$ws.Range("C21").Value = 294
$ws.Range("D21").Value = 315
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 1163
$ws.Range("G21").Value = 1293
$ws.Range("H21").Value = -10.054137664346
$ws.Range("I21").Value = 6888
$ws.Range("J21").Value = 7545
$ws.Range("K21").Value = -8.707753479125
$ws.Range("L21").Value = -10.892626131953
$ws.Range("M21").Value = -10.858030283421
$ws.Range("N21").Value = -80.563785659866
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -60
$ws.Range("G22").Value = 27
$ws.Range("H22").Value = -29.629629629629
$ws.Range("I22").Value = 93
$ws.Range("J22").Value = 110
$ws.Range("K22").Value = -15.454545454545
$ws.Range("L22").Value = 1.086956521739
$ws.Range("M22").Value = -32.116788321167
# Row 23
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 140
$ws.Range("F23").Value = 30
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 221
$ws.Range("J23").Value = 256
$ws.Range("K23").Value = -13.671875
$ws.Range("L23").Value = -12.301587301587
$ws.Range("M23").Value = 56.737588652482
# Row 24
$ws.Range("C24").Value = 249
$ws.Range("D24").Value = 325
$ws.Range("E24").Value = -23.384615384615
$ws.Range("F24").Value = 970
$ws.Range("G24").Value = 1272
$ws.Range("H24").Value = -23.742138364779
$ws.Range("I24").Value = 6721
$ws.Range("J24").Value = 7386
$ws.Range("K24").Value = -9.0035201733
$ws.Range("L24").Value = -15.288631207461
$ws.Range("M24").Value = 17.397379912663
# Row 25
$ws.Range("C25").Value = 79
$ws.Range("D25").Value = 173
$ws.Range("E25").Value = -54.335260115606
$ws.Range("F25").Value = 340
$ws.Range("G25").Value = 619
$ws.Range("H25").Value = -45.072697899838
$ws.Range("I25").Value = 2841
$ws.Range("J25").Value = 3578
$ws.Range("K25").Value = -20.598099496925
$ws.Range("L25").Value = -25.7643062451
# Row 26
$ws.Range("C26").Value = 162
$ws.Range("D26").Value = 134
$ws.Range("E26").Value = 20.895522388059
$ws.Range("F26").Value = 565
$ws.Range("G26").Value = 558
$ws.Range("H26").Value = 1.254480286738
$ws.Range("I26").Value = 3126
$ws.Range("J26").Value = 3101
$ws.Range("K26").Value = 0.806191551112
$ws.Range("L26").Value = 10.30345800988
$ws.Range("M26").Value = -5.758215254748
# Row 27
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = -58.333333333333
$ws.Range("F27").Value = 23
$ws.Range("G27").Value = 27
$ws.Range("H27").Value = -14.814814814814
$ws.Range("I27").Value = 149
$ws.Range("J27").Value = 162
$ws.Range("K27").Value = -8.024691358024
$ws.Range("L27").Value = -5.095541401273
# Row 28
$ws.Range("C28").Value = 12
$ws.Range("D28").Value = 9
$ws.Range("E28").Value = 33.333333333333
$ws.Range("G28").Value = 62
$ws.Range("H28").Value = -17.741935483871
$ws.Range("I28").Value = 318
$ws.Range("J28").Value = 353
$ws.Range("K28").Value = -9.915014164305
$ws.Range("L28").Value = 1.923076923076
# Row 29
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 25
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 13
$ws.Range("H29").Value = -7.692307692307
$ws.Range("I29").Value = 52
$ws.Range("J29").Value = 54
$ws.Range("K29").Value = -3.703703703703
$ws.Range("L29").Value = -29.729729729729
$ws.Range("M29").Value = -55.555555555555
$ws.Range("N29").Value = -86.206896551724
# Row 30
$ws.Range("C30").Value = 3
$ws.Range("E30").Value = 50
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = -10
$ws.Range("I30").Value = 45
$ws.Range("J30").Value = 49
$ws.Range("K30").Value = -8.163265306122
$ws.Range("L30").Value = -27.419354838709
$ws.Range("M30").Value = -52.631578947368
$ws.Range("N30").Value = -86.486486486486
# Row 31
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 11
$ws.Range("G31").Value = 19
$ws.Range("H31").Value = -42.105263157894
$ws.Range("I31").Value = 72
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = -10
$ws.Range("L31").Value = 71.428571428571
# Row 33
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = -50
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = -33.333333333333
$ws.Range("I33").Value = 17
$ws.Range("J33").Value = 27
$ws.Range("K33").Value = -37.037037037037
$ws.Range("L33").Value = -19.047619047619

Write-Output "edits applied"
